$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename sheets (iCC390 -> iCC389, iCC464 -> iCC470, iCC644 -> iCC651)
#    iCC431 keeps its name. Renaming also keeps definedNames ( _FilterDatabase )
#    in sync automatically since it references the sheet object.
# ---------------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("iCC390")
$wsB = $wb.Worksheets.Item("iCC431")
$wsC = $wb.Worksheets.Item("iCC464")
$wsD = $wb.Worksheets.Item("iCC644")

$wsA.Name = "iCC389"
$wsC.Name = "iCC470"
$wsD.Name = "iCC651"

# ---------------------------------------------------------------------------
# 2. sheet1 (iCC389): swap shared-string values between A15 and A37
# ---------------------------------------------------------------------------
$v15 = $wsA.Range("A15").Value2
$v37 = $wsA.Range("A37").Value2
$wsA.Range("A15").Value = $v37
$wsA.Range("A37").Value = $v15

# ---------------------------------------------------------------------------
# 3. sheet4 (iCC651): insert a new row above the former row 75 containing
#    "EX_alac__S_e" / 0, shifting the rest of the table down by one row.
# ---------------------------------------------------------------------------
$wsD.Rows.Item(75).Insert()
$wsD.Range("A75").Value = "EX_alac__S_e"
$wsD.Range("B75").Value = 0

# ---------------------------------------------------------------------------
# 4. Update the active sheet / selections to match the new view state.
#    iCC389 becomes the active tab with selection A16.
#    iCC431 keeps selection A23 (view scroll position cannot be controlled
#    from this headless engine, so we only restore the cell selection).
#    iCC470 selection becomes B11:B12 with B11 active (tab no longer selected
#    because iCC389 becomes active instead).
#    iCC651 selection becomes B3:B4 with B3 active.
# ---------------------------------------------------------------------------
$wsB.Activate() | Out-Null
$wsB.Range("A23").Select() | Out-Null

$wsC.Activate() | Out-Null
$wsC.Range("B11:B12").Select() | Out-Null

$wsD.Activate() | Out-Null
$wsD.Range("B3:B4").Select() | Out-Null

$wsA.Activate() | Out-Null
$wsA.Range("A16").Select() | Out-Null
